# "added rough title page"
#
# The deck's single slide uses the Title-Slide layout and still has its
# two placeholders (ctrTitle / subTitle) empty. Fill them in with a
# rough title page: a title and a presenter name.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$placeholders = $s.Shapes.Placeholders

# Title placeholder ("Title 1", type ctrTitle) -> "Android Presentation"
$titleRange = $placeholders.Item(1).TextFrame.TextRange
$titleRange.Text = "Android Presentation"

# Subtitle placeholder ("Subtitle 2", type subTitle) -> presenter's name.
# Typed as the first name followed by the (less common) surname, which is
# how PowerPoint naturally splits the line into two runs once the spell
# checker flags the surname.
$subtitleRange = $placeholders.Item(2).TextFrame.TextRange
$subtitleRange.Text = "Yonatan Giventer"
$lastName = $subtitleRange.Characters(9, 8)
$lastName.Text = "Giventer"
